$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- C15 / C16: "open" -> "done" ---
$ws.Range("C15").Value = "done"
$ws.Range("C16").Value = "done"

# --- Row 19: fill in the new ToDo item ---
$ws.Range("B19").Value = "Sync Noramal tabs and Customer tab "
$ws.Range("C19").Value = "open"
$ws.Range("D19").Value = "Difficult"
$ws.Range("E19").Value = "When register modified in customer, then should update other tab which contains the same register."

$ws.Rows.Item(19).RowHeight = 27

$ws.Range("B19:D19").VerticalAlignment = -4108

$ws.Range("E19").HorizontalAlignment = -4131
$ws.Range("E19").VerticalAlignment = -4108
$ws.Range("E19").WrapText = $true

# --- column widths ---
$ws.Columns.Item(4).ColumnWidth = 10.125
$ws.Columns.Item(5).ColumnWidth = 52.625

# --- selection ---
$ws.Range("B18").Select()
